# Insert a new data row for the "LHE" (Lahore, Pakistan) colo just above
# the existing "IAD" row (row 270), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 270 (pushes former row 270.. down to 271..)
$ws.Rows.Item(270).Insert()

# Copy formatting (border/bold/alignment) of the colo-code cell from the row
# that was just pushed down (now row 271, still holding the original "IAD"
# formatting) into the new blank colo-code cell so it matches the rest of
# column A.
$ws.Cells.Item(271, 1).Copy()
$ws.Cells.Item(270, 1).PasteSpecial(-4122)

# Populate the new row with the Lahore, Pakistan colo entry.
$ws.Cells.Item(270, 1).Value = "LHE"
$ws.Cells.Item(270, 2).Value = "Lahore, Pakistan"
$ws.Cells.Item(270, 3).Value = "Asia Pacific"
$ws.Cells.Item(270, 4).Value = "Lahore"
$ws.Cells.Item(270, 5).Value = "Pakistan"
$ws.Cells.Item(270, 6).Value = "PK"
$ws.Cells.Item(270, 7).Value = 31.5216007233
$ws.Cells.Item(270, 8).Value = 74.4036026001
